$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new hourly reading for 2026/01/23 needs to be inserted in date order,
# right before the existing 2026/12/29 block (row 698). Insert a row there,
# which shifts every row from 698..739 down to 699..740.
$ws.Rows.Item(698).Insert()

# Fill in the newly-opened row 698 with the new reading.
# Force column A to stay plain text ("@") so the "YYYY/MM/DD" string isn't
# auto-converted into a date serial number, matching how every other row in
# this column is stored (inline/shared text, not a date value). Resetting
# the style back to "Normal" afterwards avoids leaving a lingering custom
# number-format style on the cell.
$ws.Range("A698").NumberFormat = "@"
$ws.Range("A698").Value = "2026/01/23"
$ws.Range("A698").Style = "Normal"

$ws.Range("B698").Value = "金"
$ws.Range("C698").Value = 13
$ws.Range("D698").Value = 201
